$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 462
$ws.Range("I53").Value = 147
$ws.Range("J53").Value = 777
$ws.Range("K53").Value = 147
$ws.Range("L53").Value = 777
$ws.Range("M53").Value = 490
$ws.Range("N53").Value = -2051
$ws.Range("H62").Value = 7406.533
$ws.Range("I62").Value = 7703.5
$ws.Range("K62").Value = 7703.5
$ws.Range("M62").Value = -7079.5
$ws.Range("H65").Value = 7406.533
$ws.Range("I65").Value = 7703.5
$ws.Range("K65").Value = 38517.5
$ws.Range("M65").Value = -35397.5
$ws.Range("H132").Value = 5555.577
$ws.Range("I132").Value = 6208.4375
$ws.Range("J132").Value = 4511
$ws.Range("K132").Value = 18625.3125
$ws.Range("L132").Value = 13533
$ws.Range("M132").Value = -16095.3125
$ws.Range("N132").Value = -18593
$ws.Range("H137").Value = 5644035
$ws.Range("I137").Value = 938378.4399999999
$ws.Range("J137").Value = 9625745
$ws.Range("K137").Value = 2815135.32
$ws.Range("L137").Value = 28877235
$ws.Range("M137").Value = -2812585.32
$ws.Range("N137").Value = -28882335
$ws.Range("H138").Value = 5314.2324
$ws.Range("I138").Value = 1623.9286
$ws.Range("J138").Value = 7095.759
$ws.Range("K138").Value = 4871.7858
$ws.Range("L138").Value = 21287.277
$ws.Range("M138").Value = 268.2142000000003
$ws.Range("N138").Value = -31567.277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1892.6836
$ws.Range("I32").Value = 1814.3846
$ws.Range("K32").Value = 1814.3846
$ws.Range("M32").Value = -1527.3846
$ws.Range("H45").Value = 128571.35
$ws.Range("I45").Value = 180172.08
$ws.Range("J45").Value = 4729.6
$ws.Range("K45").Value = 180172.08
$ws.Range("L45").Value = 4729.6
$ws.Range("M45").Value = -179795.08
$ws.Range("N45").Value = -5483.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 33809.742
$ws.Range("I105").Value = 42651.074
$ws.Range("K105").Value = 42651.074
$ws.Range("M105").Value = -40904.074
$ws.Range("H134").Value = 4808.8
$ws.Range("I134").Value = 4854.7856
$ws.Range("J134").Value = 4165
$ws.Range("K134").Value = 14564.3568
$ws.Range("L134").Value = 12495
$ws.Range("M134").Value = -12029.3568
$ws.Range("N134").Value = -17565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 25641544
$ws.Range("J22").Value = 76923070
$ws.Range("L22").Value = 76923070
$ws.Range("N22").Value = -76923770
$ws.Range("H33").Value = 1666.3334
$ws.Range("I33").Value = 1666.3334
$ws.Range("K33").Value = 1666.3334
$ws.Range("M33").Value = -1287.3334
$ws.Range("H58").Value = 6387.737
$ws.Range("J58").Value = 3926
$ws.Range("L58").Value = 3926
$ws.Range("N58").Value = -4332
$ws.Range("H86").Value = 9619.6875
$ws.Range("I86").Value = 10227.875
$ws.Range("K86").Value = 10227.875
$ws.Range("M86").Value = -9104.875
$ws.Range("H88").Value = 39932.332
$ws.Range("J88").Value = 34898.5
$ws.Range("L88").Value = 34898.5
$ws.Range("N88").Value = -35710.5
$ws.Range("H89").Value = 9619.6875
$ws.Range("I89").Value = 10227.875
$ws.Range("K89").Value = 51139.375
$ws.Range("M89").Value = -45523.375
$ws.Range("H91").Value = 39932.332
$ws.Range("J91").Value = 34898.5
$ws.Range("L91").Value = 34898.5
$ws.Range("N91").Value = -37706.5
$ws.Range("H99").Value = 7356518
$ws.Range("I99").Value = 13891907
$ws.Range("J99").Value = 4205.375
$ws.Range("K99").Value = 13891907
$ws.Range("L99").Value = 4205.375
$ws.Range("M99").Value = -13890409
$ws.Range("N99").Value = -7201.375
$ws.Range("H126").Value = 7356518
$ws.Range("I126").Value = 13891907
$ws.Range("J126").Value = 4205.375
$ws.Range("K126").Value = 41675721
$ws.Range("L126").Value = 12616.125
$ws.Range("M126").Value = -41673251
$ws.Range("N126").Value = -17556.125
$ws.Range("H132").Value = 22258384
$ws.Range("I132").Value = 27789022
$ws.Range("J132").Value = 135833.33
$ws.Range("K132").Value = 83367066
$ws.Range("L132").Value = 407499.99
$ws.Range("M132").Value = -83364536
$ws.Range("N132").Value = -412559.99
$ws.Range("H134").Value = 1844219.8
$ws.Range("I134").Value = 3297564
$ws.Range("J134").Value = 3316.8667
$ws.Range("K134").Value = 9892692
$ws.Range("L134").Value = 9950.6001
$ws.Range("M134").Value = -9890157
$ws.Range("N134").Value = -15020.6001
$ws.Range("H136").Value = 6387.737
$ws.Range("J136").Value = 3926
$ws.Range("L136").Value = 11778
$ws.Range("N136").Value = -16878

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 385851.34
$ws.Range("J5").Value = 835014.5
$ws.Range("L5").Value = 2505043.5
$ws.Range("N5").Value = -2505267.5
$ws.Range("H56").Value = 4588.095
$ws.Range("I56").Value = 4588.095
$ws.Range("K56").Value = 4588.095
$ws.Range("M56").Value = -4058.095
$ws.Range("H68").Value = 3996.375
$ws.Range("J68").Value = 4447.8696
$ws.Range("L68").Value = 13343.6088
$ws.Range("N68").Value = -14965.6088
$ws.Range("H71").Value = 3996.375
$ws.Range("J71").Value = 4447.8696
$ws.Range("L71").Value = 40030.8264
$ws.Range("N71").Value = -48142.8264
$ws.Range("H107").Value = 1679.4333
$ws.Range("J107").Value = 1735.1428
$ws.Range("L107").Value = 5205.428400000001
$ws.Range("N107").Value = -9045.428400000001
$ws.Range("H132").Value = 5573058.5
$ws.Range("J132").Value = 6192188
$ws.Range("L132").Value = 55729692
$ws.Range("N132").Value = -55734752
$ws.Range("H134").Value = 12788.889
$ws.Range("I134").Value = 12788.889
$ws.Range("K134").Value = 38366.667
$ws.Range("M134").Value = -33296.667
$ws.Range("H135").Value = 385851.34
$ws.Range("J135").Value = 835014.5
$ws.Range("L135").Value = 7515130.5
$ws.Range("N135").Value = -7520200.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12845
$ws.Range("I122").Value = 13275.385
$ws.Range("K122").Value = 39826.155
$ws.Range("M122").Value = -37376.155
$ws.Range("H126").Value = 16313.777
$ws.Range("J126").Value = 11565.1
$ws.Range("L126").Value = 34695.3
$ws.Range("N126").Value = -39635.3
$ws.Range("H132").Value = 7202.579
$ws.Range("I132").Value = 5815
$ws.Range("K132").Value = 17445
$ws.Range("M132").Value = -14915
$ws.Range("H136").Value = 113189
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1426.7142
$ws.Range("I16").Value = 1622.25
$ws.Range("J16").Value = 1166
$ws.Range("K16").Value = 1622.25
$ws.Range("L16").Value = 1166
$ws.Range("M16").Value = -1452.25
$ws.Range("N16").Value = -1506
$ws.Range("H40").Value = 64214.152
$ws.Range("I40").Value = 90537.5
$ws.Range("J40").Value = 22096.8
$ws.Range("K40").Value = 90537.5
$ws.Range("L40").Value = 22096.8
$ws.Range("M40").Value = -90401.5
$ws.Range("N40").Value = -22368.8
$ws.Range("H136").Value = 5512.514
$ws.Range("I136").Value = 2146.9412
$ws.Range("K136").Value = 6440.823600000001
$ws.Range("M136").Value = -3890.823600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26613.592
$ws.Range("I132").Value = 38717.625
$ws.Range("K132").Value = 116152.875
$ws.Range("M132").Value = -113622.875
$ws.Range("H136").Value = 3932.6316
$ws.Range("I136").Value = 3159.3928
$ws.Range("J136").Value = 6097.7
$ws.Range("K136").Value = 9478.178400000001
$ws.Range("L136").Value = 18293.1
$ws.Range("M136").Value = -6928.178400000001
$ws.Range("N136").Value = -23393.1

# Row 136 on GSM loses its LeveProfitNQ (M) value entirely in this update
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M136").ClearContents()
